$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map: row number -> @(D_value, E_value) ; only price (D) and volume (E) change
# unless a coin identity swap is noted (B, C columns too).

$updates = @{
    2  = @("28.860.03", "  +1.65%  ")
    3  = @("1.882.27",  "  +0.59%  ")
    4  = @("1.004",     "  -1.27%  ")
    5  = @("314.50",    "  -0.73%  ")
    6  = @("1.004",     "  -1.38%  ")
    7  = @("0.5090",    "  -0.46%  ")
    8  = @("0.3921",    "  -1.13%  ")
    9  = @("0.08378",   "  -1.26%  ")
    10 = @("42.32",     "  +0.98%  ")
    11 = @("1.111",     "  +0.03%  ")
    12 = @("6.202",     "  -0.75%  ")
    13 = @("1.870.65",  "  +0.04%  ")
    14 = @("20.44",     "  -0.15%  ")
    15 = @("7.277",     "  +0.86%  ")
    16 = @("1.007",     "  -0.91%  ")
    17 = @("93.17",     "  +2.57%  ")
    18 = @("0.00001101","  -0.79%  ")
    19 = @("0.06712",   "  -0.81%  ")
    20 = @("17.68",     "  -0.14%  ")
    21 = @("1.005",     "  -1.20%  ")
    22 = @("5.949",     "  -0.01%  ")
    23 = @("28.797.61", "  +1.30%  ")
    24 = @("11.11",     "  -0.31%  ")
    25 = @("2.224",     "  -2.58%  ")
    26 = @("2.085.93",  "  +0.21%  ")
    29 = @("2.418",     "  +2.35%  ")
    30 = @("126.68",    "  -0.55%  ")
    33 = @("5.798",     "  +0.30%  ")
    34 = @("3.661",     "  +0.59%  ")
    35 = @("0.02454",   "  +0.73%  ")
    36 = @("0.06561",   "  +1.65%  ")
    37 = @("9.020",     "  +2.00%  ")
    38 = @("0.2164",    "  -0.69%  ")
    39 = @("5.079",     "  +1.70%  ")
    40 = @("1.199",     "  +1.15%  ")
    41 = @("1.239",     "  -1.65%  ")
    42 = @("0.6387",    "  +0.25%  ")
    43 = @("11.18",     "  -0.38%  ")
    44 = @("1.004",     "  -1.14%  ")
    45 = @("0.5997",    "  -0.22%  ")
    46 = @("13.03",     "  +0.56%  ")
    47 = @("3.675",     "  -0.84%  ")
    48 = @("2.007",     "  +0.87%  ")
    49 = @("1.226",     "  +1.85%  ")
    50 = @("122.33",    "  +0.58%  ")
    51 = @("1.168",     "  -2.92%  ")
}

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]
    $ws.Cells.Item($row, 4).Value = $vals[0]
    $ws.Cells.Item($row, 5).Value = $vals[1]
}

# Rows 27/28 swap coin identity (Monero <-> EthereumClassic) with new D/E values
$ws.Cells.Item(27, 2).Value = "EthereumClassic"
$ws.Cells.Item(27, 3).Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Cells.Item(27, 4).Value = "20.73"
$ws.Cells.Item(27, 5).Value = "  -0.15%  "

$ws.Cells.Item(28, 2).Value = "Monero"
$ws.Cells.Item(28, 3).Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Cells.Item(28, 4).Value = "157.11"
$ws.Cells.Item(28, 5).Value = "  -2.85%  "

# Rows 31/32 swap coin identity (Stellar <-> ImmutableX) with new D/E values
$ws.Cells.Item(31, 2).Value = "ImmutableX"
$ws.Cells.Item(31, 3).Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Cells.Item(31, 4).Value = "1.052"
$ws.Cells.Item(31, 5).Value = "  +1.25%  "

$ws.Cells.Item(32, 2).Value = "Stellar"
$ws.Cells.Item(32, 3).Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Cells.Item(32, 4).Value = "0.1039"
$ws.Cells.Item(32, 5).Value = "  -1.29%  "
